$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (columns D through P) with revised figures ---
$ws.Range("D2").Value = 0.25662079501451
$ws.Range("E2").Value = 0.103615919706583
$ws.Range("F2").Value = 1.2925548905567
$ws.Range("G2").Value = 1.05499243088912
$ws.Range("H2").Value = 1.58361149917828
$ws.Range("I2").Value = 2.47665412555525
$ws.Range("J2").Value = 0.0132620351176338
$ws.Range("K2").Value = 0.0492000032274128
$ws.Range("L2").Value = 0.0333926841209994
$ws.Range("M2").Value = 0.066976075339553
$ws.Range("N2").Value = 2423
$ws.Range("O2").Value = 1198
$ws.Range("P2").Value = 50.5571605447792

# --- Add new row 3: AEDB.CEA / MCP1_pg_ml_2015_rank / AsymptSympt ---
$ws.Range("A3").Value = "AEDB.CEA"
$ws.Range("B3").Value = "MCP1_pg_ml_2015_rank"
$ws.Range("C3").Value = "AsymptSympt"
$ws.Range("D3").Value = 0.333701213521244
$ws.Range("E3").Value = 0.109349512370776
$ws.Range("F3").Value = 1.39612593769705
$ws.Range("G3").Value = 1.1267937037459
$ws.Range("H3").Value = 1.72983539704799
$ws.Range("I3").Value = 3.05169365904211
$ws.Range("J3").Value = 0.0022755420659224
$ws.Range("K3").Value = 0.0531850916171955
$ws.Range("L3").Value = 0.0360284731597776
$ws.Range("M3").Value = 0.0722905320328805
$ws.Range("N3").Value = 2423
$ws.Range("O3").Value = 1199
$ws.Range("P3").Value = 50.5158893933141

# --- Add new row 4: AEDB.CEA / MCP1_rank / AsymptSympt ---
$ws.Range("A4").Value = "AEDB.CEA"
$ws.Range("B4").Value = "MCP1_rank"
$ws.Range("C4").Value = "AsymptSympt"
$ws.Range("D4").Value = 0.365453270635651
$ws.Range("E4").Value = 0.125363018444617
$ws.Range("F4").Value = 1.44116709885159
$ws.Range("G4").Value = 1.12720571816695
$ws.Range("H4").Value = 1.84257635792503
$ws.Range("I4").Value = 2.91516010997375
$ws.Range("J4").Value = 0.00355506107610972
$ws.Range("K4").Value = 0.0327959127873894
$ws.Range("L4").Value = 0.0278573516612683
$ws.Range("M4").Value = 0.0482411855261864
$ws.Range("N4").Value = 2423
$ws.Range("O4").Value = 556
$ws.Range("P4").Value = 77.05323978539
